# Insert a new data colo row ("LHE" / Lahore, Pakistan) at sheet row 270.
# This pushes the existing rows 270..328 down to 271..329 (e.g. IAD now at
# 271, and the final Halifax/YHZ row moves from 328 to 329), and grows the
# sheet dimension from A1:H328 to A1:H329.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 270 (the row that currently
# holds IAD / Ashburn), shifting everything below it down by one.
$ws.Rows(270).Insert()

# The newly inserted row is blank/unformatted. Copy the formatting
# (bold + border + center/top alignment on column A) from the row that is
# now directly below it (the old IAD row, now at 271) so the new row
# matches the look of every other data row.
$ws.Range("A271").Copy()
$ws.Range("A270").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Populate the new row with the Lahore, Pakistan colo entry.
$ws.Range("A270").Value = "LHE"
$ws.Range("B270").Value = "Lahore, Pakistan"
$ws.Range("C270").Value = "Asia Pacific"
$ws.Range("D270").Value = "Lahore"
$ws.Range("E270").Value = "Pakistan"
$ws.Range("F270").Value = "PK"
$ws.Range("G270").Value = 31.5216007233
$ws.Range("H270").Value = 74.4036026001
